# New pressure sensor added to the controller
# 1) Append a new measurement row (index 19) to Table1, and
# 2) relocate the scratch/"note" calculations block from H27:W45 to E27:T45
#    (3 columns to the left), matching the author's rearrangement.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Part 1: add the new row (Index 19) to Table1, which grows B3:L21 -> B3:L22
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$null = $lo.ListRows.Add()

# Bring over number formats (date/time columns) from the row above so the
# new row matches the rest of the table's look.
$ws.Range("G21:H21").Copy()
$ws.Range("G22").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("O21:S21").Copy()
$ws.Range("O22").PasteSpecial(-4122)
$ws.Range("U21:V21").Copy()
$ws.Range("U22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B22").Value = 19
$ws.Range("C22").Value = 286
$ws.Range("D22").Value = 159
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 45887.504166666666
$ws.Range("H22").Value = 45888.709722222222
$ws.Range("I22").Value = 11.9
$ws.Range("J22").Value = 9
$ws.Range("K22").Value = 2.5

$ws.Range("O22").Formula = "=H22-G22"
$ws.Range("P22").Formula = "=O22"
$ws.Range("Q22").Formula = "=I22-J22"
$ws.Range("R22").Formula = "=(F22-E22)/0.9982"
$ws.Range("S22").Formula = "=K22*P22"
$ws.Range("U22").Formula = "=Q22*1440/1736"
$ws.Range("V22").Formula = "=(1-ABS(U22-K22)/K22)*100%"

# Extend the color-scale conditional formatting that covers the Accuracy
# column so it still spans the whole table (V4:V21 -> V4:V22).
$ws.Range("V4:V21").FormatConditions.Delete()
$null = $ws.Range("V4:V22").FormatConditions.AddColorScale(3)

# ---------------------------------------------------------------------
# Part 2: move the scratch calculations block 3 columns to the left
# (H27:W45 -> E27:T45). Recreated cell-by-cell (values/formulas/styles)
# since these rows/columns sit well outside any ListObject.
# ---------------------------------------------------------------------

# Preserve the custom number format used on the two "raw" weight cells
# before the source range is cleared.
$ws.Range("H28:H29").Copy()
$ws.Range("E28").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("H27:W45").Clear()

# Row 27
$ws.Range("I27").Value = 65.77
$ws.Range("M27").Value = 80.89
$ws.Range("Q27").Value = 74

# Row 28
$ws.Range("E28").Value = 65.773633770000004
$ws.Range("F28").Value = 2.3
$ws.Range("G28").Formula = "=E28/F28"
$ws.Range("I28").Value = 247
$ws.Range("J28").Value = 0.4
$ws.Range("K28").Value = 0.43
$ws.Range("L28").Formula = "=J28-K28"
$ws.Range("M28").Value = 1367
$ws.Range("N28").Value = 2.4
$ws.Range("O28").Value = 2.37
$ws.Range("P28").Formula = "=N28-O28"
$ws.Range("Q28").Value = 940
$ws.Range("R28").Value = 1.55
$ws.Range("S28").Value = 1.63
$ws.Range("T28").Formula = "=R28-S28"

# Row 29
$ws.Range("E29").Value = 80.890960000000007
$ws.Range("F29").Value = 2.69
$ws.Range("G29").Formula = "=E29/F29"
$ws.Range("I29").Value = 1007
$ws.Range("J29").Value = 1.75
$ws.Range("K29").Value = 1.75
$ws.Range("L29").Formula = "=J29-K29"
$ws.Range("M29").Value = 4312
$ws.Range("N29").Value = 8.05
$ws.Range("O29").Value = 7.49
$ws.Range("P29").Formula = "=N29-O29"
$ws.Range("Q29").Value = 2516
$ws.Range("R29").Value = 4.15
$ws.Range("S29").Value = 4.37
$ws.Range("T29").Formula = "=R29-S29"

# Row 30
$ws.Range("F30").Value = 2.5
$ws.Range("G30").Formula = "=(G28+G29)/2"
$ws.Range("I30").Value = 1063
$ws.Range("J30").Value = 1.65
$ws.Range("K30").Value = 1.85
$ws.Range("L30").Formula = "=J30-K30"
$ws.Range("M30").Value = 5679
$ws.Range("N30").Value = 10.45
$ws.Range("O30").Value = 9.86
$ws.Range("P30").Formula = "=N30-O30"

# Row 31
$ws.Range("G31").Formula = "=G30*F30"

# Row 33
$ws.Range("I33").Formula = "=247/-0.03"
$ws.Range("M33").Formula = "=M28/P28"
$ws.Range("Q33").Formula = "=Q28/T28"

# Row 34
$ws.Range("I34").Formula = "=1063/-0.2"
$ws.Range("M34").Formula = "=M29/P29"
$ws.Range("Q34").Formula = "=Q29/T29"

# Row 35
$ws.Range("M35").Formula = "=M30/P30"

# Row 37
$ws.Range("I37").Formula = "=5679/1063"
$ws.Range("Q37").Formula = "=M30/Q29"

# Row 38
$ws.Range("I38").Formula = "=I30*I37"
$ws.Range("L38").Formula = "=L30*I37"
$ws.Range("Q38").Formula = "=Q29*Q37"
$ws.Range("T38").Formula = "=T29*Q37"

# Row 44
$ws.Range("I44").Value = 65.77
$ws.Range("K44").Value = 80.89
$ws.Range("M44").Value = 74
$ws.Range("O44").Value = 77.150000000000006

# Row 45
$ws.Range("I45").Formula = "=L38"
$ws.Range("K45").Formula = "=P30"
$ws.Range("M45").Formula = "=T38"
$ws.Range("O45").Value = 0

# ---------------------------------------------------------------------
# View state: cursor ends up where the author left it.
# ---------------------------------------------------------------------
$ws.Range("P52").Select()
